# Auto-generated: applies latest crypto price/volume snapshot to sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row=2; D='66.010.21'; E='  -0.03%  ' },
    @{ Row=3; D='3.328.30'; E='  +1.46%  ' },
    @{ Row=4; D='1.00'; E='  +0.03%  ' },
    @{ Row=5; B='BNB'; C='https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb'; D='558.86'; E='  +0.14%  ' },
    @{ Row=6; B='Solana'; C='https://coinranking.com/coin/zNZHO_Sjf+solana-sol'; D='186.08'; E='  +0.34%  ' },
    @{ Row=7; D='1.00'; E='  +0.07%  ' },
    @{ Row=8; D='3.319.84'; E='  +1.40%  ' },
    @{ Row=9; D='0.577'; E='  -2.66%  ' },
    @{ Row=10; D='0.176'; E='  -6.39%  ' },
    @{ Row=11; D='0.578'; E='  -1.51%  ' },
    @{ Row=12; D='45.89'; E='  -3.28%  ' },
    @{ Row=13; D='0.0000263'; E='  -1.51%  ' },
    @{ Row=14; D='3.862.06'; E='  +1.50%  ' },
    @{ Row=15; D='8.46'; E='  -1.58%  ' },
    @{ Row=16; D='571.55'; E='  -9.51%  ' },
    @{ Row=17; D='66.032.56' },
    @{ Row=18; B='WrappedEther'; C='https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'; D='3.333.66'; E='  +1.69%  ' },
    @{ Row=19; B='TRON'; C='https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'; D='0.117'; E='  +0.58%  ' },
    @{ Row=20; D='17.71'; E='  -1.04%  ' },
    @{ Row=21; D='10.88'; E='  -4.37%  ' },
    @{ Row=22; E='  -1.32%  ' },
    @{ Row=23; D='18.13'; E='  +1.45%  ' },
    @{ Row=24; E='  +1.30%  ' },
    @{ Row=25; D='97.64'; E='  -8.46%  ' },
    @{ Row=26; D='3.96'; E='  -0.54%  ' },
    @{ Row=27; D='2.70'; E='  +0.92%  ' },
    @{ Row=28; D='9.40'; E='  -1.64%  ' },
    @{ Row=29; D='8.45'; E='  -3.00%  ' },
    @{ Row=30; D='30.60'; E='  +0.50%  ' },
    @{ Row=31; D='6.69'; E='  +6.67%  ' },
    @{ Row=32; B='Bittensor'; C='https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'; D='572.00'; E='  +5.83%  ' },
    @{ Row=33; B='dogwifhat'; C='https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'; D='3.69'; E='  -9.22%  ' },
    @{ Row=34; E='  -1.78%  ' },
    @{ Row=35; E='  -1.73%  ' },
    @{ Row=36; D='3.731.01'; E='  +1.64%  ' },
    @{ Row=37; E='  +0.09%  ' },
    @{ Row=38; D='55.55'; E='  -3.17%  ' },
    @{ Row=39; D='33.97'; E='  +4.24%  ' },
    @{ Row=40; D='0.126'; E='  -3.59%  ' },
    @{ Row=41; D='0.0₃0690'; E='  -6.14%  ' },
    @{ Row=42; D='2.60'; E='  -5.02%  ' },
    @{ Row=43; D='3.12'; E='  -8.90%  ' },
    @{ Row=44; D='3.37'; E='  +2.92%  ' },
    @{ Row=45; E='  -0.98%  ' },
    @{ Row=46; E='  -1.51%  ' },
    @{ Row=47; B='CoreDAO'; C='https://coinranking.com/coin/HFvoXUQh4+coredao-core'; D='2.96'; E='  -11.15%  ' },
    @{ Row=48; B='FirstDigitalUSD'; C='https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'; D='0.999'; E='  +0.17%  ' },
    @{ Row=49; B='Stellar'; C='https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'; D='0.126'; E='  -2.43%  ' },
    @{ Row=50; D='2.51'; E='  -4.00%  ' },
    @{ Row=51; D='125.46'; E='  +2.59%  ' }
)

foreach ($u in $updates) {
    $r = $u.Row
    foreach ($col in @('B','C','D','E')) {
        if ($u.ContainsKey($col)) {
            $ref = "$col$r"
            $cell = $ws.Range($ref)
            $cell.NumberFormat = "@"
            $cell.Value = $u[$col]
        }
    }
}
